$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the existing "_GoBack" bookmark (it currently sits at the
#    end of the "An H-L quadrant..." paragraph). It will be re-added
#    later at the new insertion point.
# ------------------------------------------------------------------
try {
    $oldBm = $d.Bookmarks.Item("_GoBack")
    $oldBm.Delete()
} catch {
    # no pre-existing _GoBack bookmark - nothing to remove
}

# ------------------------------------------------------------------
# 2) Append new sentences to the paragraph ending in
#    "... I plot was slightly different. " and drop the new
#    "_GoBack" bookmark right after the new text.
# ------------------------------------------------------------------
$para = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "*I plot was slightly different.*") {
        $para = $candidate
    }
}
if ($null -eq $para) {
    throw "Could not find anchor paragraph text"
}

$pStart = $para.Range.Start
$pEnd = $para.Range.End - 1   # exclude the paragraph mark
$target = $d.Range($pStart, $pEnd)

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r w:rsidRPr="00896772"><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="4472C4" w:themeColor="accent5"/></w:rPr><w:t xml:space="preserve">From my lab, the results were not too different from each other. The p-values were similar however the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00896772"><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="4472C4" w:themeColor="accent5"/></w:rPr><w:t>moran’s</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00896772"><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="4472C4" w:themeColor="accent5"/></w:rPr><w:t xml:space="preserve"> I plot was slightly different. </w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="4472C4" w:themeColor="accent5"/></w:rPr><w:t>This being said, the histograms p</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="4472C4" w:themeColor="accent5"/></w:rPr><w:t xml:space="preserve">roduced showed large variations, these mainly being the majority of the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="4472C4" w:themeColor="accent5"/></w:rPr><w:t>neighbours</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="4472C4" w:themeColor="accent5"/></w:rPr><w:t xml:space="preserve"> being past the half way point.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$null = $target.InsertXML($xml)
